# poisson_naive versao media ponderada
# Update column A (row index / weight) values on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 26
    3  = 28
    4  = 30
    5  = 32
    6  = 33
    7  = 36
    8  = 37
    9  = 39
    10 = 41
    11 = 44
    12 = 46
    13 = 48
    14 = 50
    15 = 14
    16 = 55
    17 = 83
    18 = 146
    19 = 164
    20 = 207
    21 = 217
    22 = 261
    23 = 333
    24 = 350
    25 = 409
    26 = 485
    27 = 516
}

foreach ($row in $newValues.Keys) {
    $ws.Range("A$row").Value = $newValues[$row]
}
